$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(5)

$tf = $sh.TextFrame
$tr = $tf.TextRange

# 1) "Menu :{" -> "Menu :"
$paraMenu = $tr.Paragraphs(4, 1)
$paraMenu.Text = "Menu :"

# 2) Demote Starter/Main/Dessert to level 2 (lvl="1" in OOXML, 0-based)
$paraStarter = $tr.Paragraphs(5, 1)
$paraStarter.IndentLevel = 2

$paraMain = $tr.Paragraphs(6, 1)
$paraMain.IndentLevel = 2

$paraDessert = $tr.Paragraphs(7, 1)
$paraDessert.IndentLevel = 2
# "Dessert:}" -> "Dessert:"
$paraDessert.Text = "Dessert:"

# 3) Remove the trailing "Guests [users ID,...]" paragraph entirely (it was
#    the last paragraph in the text body, right after "Images :"). Deleting
#    the last paragraph once only clears its text and leaves a stray empty
#    paragraph behind, so we delete it a second time to fully collapse it
#    away and restore "Images :" as the final paragraph.
$paraGuests = $tr.Paragraphs(16, 1)
$paraGuests.Delete()
$paraGuestsLeftover = $tr.Paragraphs(16, 1)
$paraGuestsLeftover.Delete()

# The shape has <a:spAutoFit/>, so its height (cy) automatically recomputes
# from the now-shorter text (fewer paragraphs) - matching the diff's
# cy="4801314" -> cy="4524315" without needing to set it explicitly.
